$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Trening" header in column F, copying the header formatting
# (bold font, border, centered alignment) from the neighboring header cell.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value2 = "Trening"

# Replace the text timestamps in column A (rows 2-7) with numeric date
# serials, and extend the data with new rows 6-13 (rows below row 5 shift
# down and new measurements are appended).
$timestamps = @(
    45685.65130983797,
    45685.66109571759,
    45685.66112465278,
    45685.65130636574,
    45685.65587581018,
    45685.66109456019,
    45685.66948229167,
    45685.67614664352,
    45685.67821261574,
    45685.66948113426,
    45685.6761443287,
    45685.67821030092
)

$seconds = @(1388.1, 2233.6, 2236.1, 1387.8, 1782.6, 2233.5, 2958.2, 3534, 3712.5, 2958.1, 3533.8, 3712.3)
$velocity = @(13.59, 11.33, 13.86, 9.720000000000001, 9.23, 9.76, 10.31, 11.41, 11.76, 9.27, 8.99, 9.52)
$accel = @(3.676808357238768, 3.210851396833146, 3.511942522866386, 3.443315301622663, 2.725719043186733, 3.048933403832571, 2.829279439789907, 3.119807311466759, 3.044049058641705, 2.919103758675711, 2.961176872253416, 2.861920424870082)
$bin = @("10-15", "10-15", "10-15", "5-10", "5-10", "5-10", "10-15", "10-15", "10-15", "5-10", "5-10", "5-10")
$trening = @("Duża Gra", "Duża Gra", "Duża Gra", "Duża Gra", "Duża Gra", "Duża Gra", "Mała Gra", "Mała Gra", "Mała Gra", "Mała Gra", "Mała Gra", "Mała Gra")

# Apply the lowercase date/time format first - this registers it as
# numFmtId 164, matching how the workbook's style table ends up looking.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value2 = $seconds[$i]
    $ws.Cells.Item($row, 3).Value2 = $velocity[$i]
    $ws.Cells.Item($row, 4).Value2 = $accel[$i]
    $ws.Cells.Item($row, 5).Value2 = $bin[$i]
    $ws.Cells.Item($row, 6).Value2 = $trening[$i]
}

# Now apply the final (uppercase) date/time format on all the timestamp
# cells - this registers numFmtId 165 and is the one actually referenced
# by the cell styles.
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
